# Auto-generated edit script applying cryptos.xlsx price/volume refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.817.30"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "'3.172.64"
$ws.Range("E3").Value = "  -4.74%  "

$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "'572.95"
$ws.Range("E5").Value = "  -0.56%  "

$ws.Range("D6").Value = "'172.21"
$ws.Range("E6").Value = "  -3.75%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").Value = "'0.596"
$ws.Range("E8").Value = "  -2.76%  "

$ws.Range("D9").Value = "'3.169.24"
$ws.Range("E9").Value = "  -4.79%  "

$ws.Range("D10").Value = "'0.125"
$ws.Range("E10").Value = "  -2.92%  "

$ws.Range("D11").Value = "'6.61"
$ws.Range("E11").Value = "  -3.70%  "

$ws.Range("D12").Value = "'0.393"
$ws.Range("E12").Value = "  -3.54%  "

$ws.Range("D13").Value = "'3.720.06"
$ws.Range("E13").Value = "  -4.89%  "

$ws.Range("D14").Value = "'0.136"
$ws.Range("E14").Value = "  +0.80%  "

$ws.Range("D15").Value = "'27.52"
$ws.Range("E15").Value = "  -4.18%  "

$ws.Range("D16").Value = "'65.762.33"
$ws.Range("E16").Value = "  +0.44%  "

$ws.Range("D17").Value = "'0.0000165"
$ws.Range("E17").Value = "  -2.31%  "

$ws.Range("D18").Value = "'3.171.18"
$ws.Range("E18").Value = "  -4.82%  "

$ws.Range("D19").Value = "'5.74"
$ws.Range("E19").Value = "  +0.16%  "

$ws.Range("D20").Value = "'12.94"
$ws.Range("E20").Value = "  -3.60%  "

$ws.Range("D21").Value = "'361.80"
$ws.Range("E21").Value = "  -0.28%  "

$ws.Range("D22").Value = "'7.30"
$ws.Range("E22").Value = "  -1.75%  "

$ws.Range("E23").Value = "  +0.22%  "

$ws.Range("D24").Value = "'69.21"
$ws.Range("E24").Value = "  -2.97%  "

$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").Value = "'0.496"
$ws.Range("E25").Value = "  -4.55%  "

$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D26").Value = "'3.307.32"
$ws.Range("E26").Value = "  -4.96%  "

$ws.Range("E27").Value = "  -6.13%  "

$ws.Range("D28").Value = "'9.90"
$ws.Range("E28").Value = "  +3.64%  "

$ws.Range("E29").Value = "  -0.63%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "'1.93"
$ws.Range("E31").Value = "  -1.50%  "

$ws.Range("B32").Value = "USDe"
$ws.Range("C32").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D32").Value = "'0.998"
$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'5.40"
$ws.Range("E33").Value = "  -3.70%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "'22.15"
$ws.Range("E34").Value = "  -3.30%  "

$ws.Range("D35").Value = "'1.21"
$ws.Range("E35").Value = "  -0.63%  "

$ws.Range("D36").Value = "'6.65"
$ws.Range("E36").Value = "  -2.35%  "

$ws.Range("D37").Value = "'160.14"
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("D38").Value = "'1.46"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("D39").Value = "'0.836"
$ws.Range("E39").Value = "  -0.94%  "

$ws.Range("E40").Value = "  +3.99%  "

$ws.Range("D41").Value = "'26.46"
$ws.Range("E41").Value = "  -3.41%  "

$ws.Range("D42").Value = "'2.50"
$ws.Range("E42").Value = "  -2.94%  "

$ws.Range("D43").Value = "'2.648.20"
$ws.Range("E43").Value = "  -2.56%  "

$ws.Range("D44").Value = "'6.18"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("D45").Value = "'4.20"
$ws.Range("E45").Value = "  -1.57%  "

$ws.Range("D46").Value = "'39.74"
$ws.Range("E46").Value = "  -0.02%  "

$ws.Range("B47").Value = "Bittensor"
$ws.Range("C47").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D47").Value = "'331.41"
$ws.Range("E47").Value = "  -1.49%  "

$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").Value = "'0.0657"
$ws.Range("E48").Value = "  -1.49%  "

$ws.Range("D49").Value = "'24.11"
$ws.Range("E49").Value = "  +0.69%  "

$ws.Range("D50").Value = "'0.0275"
$ws.Range("E50").Value = "  -1.28%  "

$ws.Range("D51").Value = "'0.102"
$ws.Range("E51").Value = "  -1.20%  "
